$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename header cells for the two columns (T, U) that were previously
# "MeanStageFirstRep" / "MeanStageFirstRepObs"
$ws.Range("T1").Value = "StageFirstRep"
$ws.Range("U1").Value = "StageFirstRepObs"

# Add "Other" to the GrowthTransition options list (column AD), inserted
# right before "Unknown" so "Unknown" shifts down one row.
$ws.Range("AD9").Insert()
$ws.Range("AD9").Value = "Other"

# Update the view state to match where the user ended up looking.
$ws.Activate()
$ws.Range("AD11").Select()
